# Heat Distribution Losses.xlsx
# "add 2019 to files, update calibration"
#  - HDL sheet: insert a new first data column (2019) ahead of 2020, so the
#    table reads 2019..2050 instead of 2020..2050; the 2019 column simply
#    mirrors the (still About!A37-calibrated) 2020 value.
#  - HDL sheet becomes the active / visible tab, with B3 selected.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$hdl = $wb.Worksheets.Item("HDL")

# --- HDL sheet: insert a new column B (year 2019), shifting 2020..2050 (old B:AF) to C:AG ---
$hdl.Columns("B:B").Insert()

# Column B should look like column A (same custom width)
$hdl.Columns("B:B").ColumnWidth = $hdl.Columns("A:A").ColumnWidth

# --- Row 1 (year headers): B1 = 2019, formatted like the "Year" header cell A1 ---
$hdl.Range("A1").Copy()
$hdl.Range("B1").PasteSpecial(-4122)
$hdl.Range("B1").Value = 2019

# --- Row 2 (calibration values) ---
# B2 mirrors the (new first) data column C2
$hdl.Range("B2").Formula = "=C2"
$hdl.Range("B2").Style = "Normal"

# C2 keeps pulling the calibrated value from the About sheet, with plain/default formatting
$hdl.Range("C2").Formula = "=About!A37"
$hdl.Range("C2").Style = "Normal"

# D2 through AG2 all reference C2 (absolute column, relative row)
$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG")
foreach ($c in $cols) {
    $hdl.Range("$c`2").Formula = "=`$C2"
}

# --- Make HDL the active/visible tab with B3 selected (About loses tabSelected) ---
$hdl.Activate()
$hdl.Range("B3").Select()

$wb.Save()
